$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TVN sampling method data (denser sampling interval) replacing the
# previous snippet rows for A2:E51.
$data = @(
    @(91, "6:10 AM", 2710.555, 2740.555, 4.655),
    @(105, "6:17 AM", 3136.858, 3166.858, 6.486),
    @(116, "6:23 AM", 3473.17, 3503.17, 1.063333),
    @(129, "6:30 AM", 3867.07, 3897.07, 0.835),
    @(145, "6:37 AM", 4332.1125, 4362.1125, 1.4225),
    @(174, "6:52 AM", 5193.04, 5223.04, 71.97),
    @(214, "7:12 AM", 6404.45, 6434.45, 1.525),
    @(251, "7:30 AM", 7516.66, 7546.66, 27.63),
    @(331, "8:10 AM", 9903.9, 9933.9, 66.37),
    @(344, "8:17 AM", 10303.675, 10333.675, 0.5),
    @(483, "9:27 AM", 14487.21, 14517.21, 5.45),
    @(509, "9:39 AM", 15258.83, 15288.83, 2.655),
    @(531, "9:50 AM", 15908.795, 15938.795, 4.8),
    @(981, "1:35 PM", 29415.15, 29445.15, 13.746667),
    @(995, "1:42 PM", 29826.6, 29856.6, 6.235),
    @(1003, "1:46 PM", 30076.196667, 30106.196667, 3.146667),
    @(1016, "1:53 PM", 30454.73, 30484.73, 24.88),
    @(1025, "1:57 PM", 30728.4775, 30758.4775, 4.9925),
    @(1031, "2:00 PM", 30916.155, 30946.155, 5.525),
    @(1037, "2:03 PM", 31099.75, 31129.75, 6.075),
    @(1052, "2:11 PM", 31541.535, 31571.535, 0.8425),
    @(1066, "2:18 PM", 31970.67, 32000.67, 17.12),
    @(1073, "2:21 PM", 32172.445, 32202.445, 12.6),
    @(1101, "2:35 PM", 33002.63, 33032.63, 2.815),
    @(1114, "2:42 PM", 33419.09, 33449.09, 2.96),
    @(1157, "3:03 PM", 34686.29, 34716.29, 3),
    @(1166, "3:08 PM", 34958.575, 34988.575, 2.23),
    @(1185, "3:17 PM", 35530.52, 35560.52, 115.925),
    @(1192, "3:21 PM", 35742.3425, 35772.3425, 13.795),
    @(1204, "3:27 PM", 36097.916667, 36127.916667, 82.06),
    @(1212, "3:31 PM", 36344.033333, 36374.033333, 13.163333),
    @(1230, "3:40 PM", 36884.13, 36914.13, 3.52),
    @(1245, "3:47 PM", 37332.09, 37362.09, 193.12),
    @(1274, "4:02 PM", 38209.85, 38239.85, 4.86),
    @(1323, "4:26 PM", 39677.255, 39707.255, 57.735),
    @(1352, "4:41 PM", 40545.91, 40575.91, 3.843333),
    @(1367, "4:48 PM", 40993.815, 41023.815, 30.44),
    @(1379, "4:54 PM", 41352.525, 41382.525, 1.775),
    @(1393, "5:01 PM", 41774.303333, 41804.303333, 14.303333),
    @(1406, "5:08 PM", 42163.55, 42193.55, 14.93),
    @(1414, "5:12 PM", 42400.085, 42430.085, 51.73),
    @(1425, "5:17 PM", 42732.5175, 42762.5175, 0.9325),
    @(1436, "5:23 PM", 43062.525, 43092.525, 11.06),
    @(1443, "5:26 PM", 43275.036667, 43305.036667, 14.33),
    @(1451, "5:30 PM", 43516.02, 43546.02, 88.58),
    @(1481, "5:45 PM", 44405.94, 44435.94, 5.48),
    @(1535, "6:12 PM", 46024.31, 46054.31, 7.985),
    @(1545, "6:17 PM", 46329.05, 46359.05, 3.473333),
    @(1554, "6:22 PM", 46615.0725, 46645.0725, 1.5575),
    @(1611, "6:50 PM", 48310.16, 48340.16, 3.305)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
